# Replace the version number "13" with "20" in the second paragraph,
# and make sure the "_GoBack" bookmark ends up wrapping the replaced
# run (bookmarkStart right before it, bookmarkEnd right after it) --
# mirroring how Word re-anchors the last-edit bookmark around an
# in-place text replacement.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("13", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "20", 2)

if ($found) {
    $editStart = $rng.Start
    $editEnd = $rng.End

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $goBackRange = $d.Range($editStart, $editEnd)
    $d.Bookmarks.Add("_GoBack", $goBackRange)
}
